# Slides for section 3
# Slide 4 ("Using Speech") title is retitled to "Predicting Taxi Journey Times".
# The original author typed/edited the text so that it ended up as two runs:
#   "Predicting Taxi " + "Journey Times"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(4)
$sh = $s.Shapes.Item(1)          # ctrTitle placeholder ("Shape 155")

$tr = $sh.TextFrame.TextRange
$tr.Text = "Predicting Taxi "
[void]$tr.InsertAfter("Journey Times")
